# Updates the "revision" (column B) values produced by the yoy GVA
# revision series recomputation (ifo GDP component analysis preprocessing).
# Only the numeric values change; dates (column A) and row 12 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 0.2031110624549193
$ws.Range("B4").Value  = -1.253592906619461
$ws.Range("B5").Value  = -0.4722879045915285
$ws.Range("B6").Value  = 0.3841888758731482
$ws.Range("B7").Value  = 0.6238557275990209
$ws.Range("B8").Value  = -0.5986487698215637
$ws.Range("B9").Value  = -0.8403875287742224
$ws.Range("B10").Value = 0.2589026302891151
$ws.Range("B11").Value = -0.01037278398443586
$ws.Range("B13").Value = 0.1028306007391011
$ws.Range("B14").Value = -0.4746561949739103
$ws.Range("B15").Value = -0.6319729579588484
$ws.Range("B16").Value = 0.2398729025141266
$ws.Range("B17").Value = -0.2240650062994343
$ws.Range("B18").Value = -0.9183314158605516
$ws.Range("B19").Value = -1.10539566971255
$ws.Range("B20").Value = -0.3764554607521564
$ws.Range("B21").Value = 0.4735660008982068
$ws.Range("B22").Value = 0.3026860865385061
